$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week#5")

$ws.Range("D2").Value = "1) Design UI for data selection and display various charts using HTML5 and JavaScript`n2) Design model"
$ws.Range("D3").Value = "1) Query Data Set to get proper Data using dataframe.js`n2) Design Model"
$ws.Range("D4").Value = "1) Design Classes in Java Script using E6`n2) Sequesnce diagram"
$ws.Range("D5").Value = "1) Use chart.js to generate bar graphs - Running example in legacy code. Understand all fields and their functionality`n2) Domain model"
$ws.Range("D6").Value = "1) Use chart.js to generate pie charts - Running example in legacy code. Understand all fields and their functionality`n2)Domain model"

$ws.Activate()
$ws.Range("F5").Select() | Out-Null
